# Logged Week 16 and performed season sim from Week 17
# New players added to the roster: M.Ffrench and J.Moore. Insert a column
# for each (shifting everything at/after the insertion point to the right)
# on both the Rushing and Receiving sheets, fill in the header name and the
# "n" placeholder used for the rest of the "Yards list" row.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Insert new column before the "K.Allen" column (directly after G.Nabers)
    # and give it the header "M.Ffrench", matching the bold/centered/bordered
    # header style used by the rest of row 1.
    $ws.Columns("J").Insert()
    $hdr = $ws.Cells.Item(1, 10)
    $hdr.Value = "M.Ffrench"
    $hdr.Font.Bold = $true
    $hdr.HorizontalAlignment = -4108
    $hdr.VerticalAlignment = -4160
    $hdr.Borders.LineStyle = 1
    $ws.Cells.Item(2, 10).Value = "n"

    # Insert new column before the "J.Cook" column (directly after A.Roberts,
    # which has by now shifted one column to the right, to column P).
    $ws.Columns("Q").Insert()
    $hdr2 = $ws.Cells.Item(1, 17)
    $hdr2.Value = "J.Moore"
    $hdr2.Font.Bold = $true
    $hdr2.HorizontalAlignment = -4108
    $hdr2.VerticalAlignment = -4160
    $hdr2.Borders.LineStyle = 1
    $ws.Cells.Item(2, 17).Value = "n"
}
